$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.051.29"
$ws.Range("E2").Value = "  -5.91%  "

$ws.Range("D3").Value = "2.447.39"
$ws.Range("E3").Value = "  -8.55%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'540.52"
$ws.Range("E5").Value = "  -2.29%  "

$ws.Range("D6").Value = "'147.29"
$ws.Range("E6").Value = "  -6.80%  "

$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("D8").Value = "'0.571"
$ws.Range("E8").Value = "  -3.31%  "

$ws.Range("D9").Value = "2.463.08"
$ws.Range("E9").Value = "  -8.02%  "

$ws.Range("D10").Value = "'0.0992"
$ws.Range("E10").Value = "  -6.32%  "

$ws.Range("E11").Value = "  -2.19%  "

$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("E13").Value = "  -4.22%  "

$ws.Range("E14").Value = "  -8.42%  "

$ws.Range("D15").Value = "'23.96"
$ws.Range("E15").Value = "  -9.32%  "

$ws.Range("D16").Value = "58.928.58"
$ws.Range("E16").Value = "  -5.98%  "

$ws.Range("E17").Value = "  -6.15%  "

$ws.Range("D18").Value = "2.516.66"
$ws.Range("E18").Value = "  -5.99%  "

$ws.Range("D19").Value = "'11.11"
$ws.Range("E19").Value = "  -6.32%  "

$ws.Range("D20").Value = "'4.36"
$ws.Range("E20").Value = "  -5.66%  "

$ws.Range("D21").Value = "'324.24"
$ws.Range("E21").Value = "  -5.85%  "

$ws.Range("D22").Value = "'0.965"
$ws.Range("E22").Value = "  -3.39%  "

$ws.Range("D23").Value = "'5.71"
$ws.Range("E23").Value = "  -9.01%  "

$ws.Range("D24").Value = "'60.72"
$ws.Range("E24").Value = "  -3.90%  "

$ws.Range("D25").Value = "'0.450"
$ws.Range("E25").Value = "  -11.37%  "

$ws.Range("E26").Value = "  -4.84%  "

$ws.Range("E27").Value = "  -2.05%  "

$ws.Range("E28").Value = "  -6.09%  "

$ws.Range("E29").Value = "  -5.82%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.26"
$ws.Range("E30").Value = "  -11.43%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0770"
$ws.Range("E31").Value = "  -9.71%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.67"
$ws.Range("E32").Value = "  -7.60%  "

$ws.Range("D33").Value = "'0.997"
$ws.Range("E33").Value = "  -0.14%  "

$ws.Range("D34").Value = "'156.63"
$ws.Range("E34").Value = "  -3.73%  "

$ws.Range("D35").Value = "'1.38"
$ws.Range("E35").Value = "  -6.65%  "

$ws.Range("E36").Value = "  -5.34%  "

$ws.Range("E37").Value = "  -9.11%  "

$ws.Range("E38").Value = "  -3.51%  "

$ws.Range("D39").Value = "'316.71"
$ws.Range("E39").Value = "  -9.71%  "

$ws.Range("E40").Value = "  -5.58%  "

$ws.Range("D41").Value = "'0.839"
$ws.Range("E41").Value = "  -11.68%  "

$ws.Range("D42").Value = "'36.22"
$ws.Range("E42").Value = "  -5.54%  "

$ws.Range("E43").Value = "  -6.75%  "

$ws.Range("E46").Value = "  -2.79%  "

$ws.Range("E47").Value = "  -5.82%  "

$ws.Range("E48").Value = "  -5.95%  "

$ws.Range("E49").Value = "  -4.97%  "

$ws.Range("D50").Value = "'121.91"
$ws.Range("E50").Value = "  -5.47%  "

$ws.Range("D51").Value = "'18.89"
$ws.Range("E51").Value = "  -9.67%  "
